$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "python_vs_r"

$black = [System.Drawing.ColorTranslator]::ToOle([System.Drawing.Color]::FromArgb(0,0,0))

$tmp1 = $ws.Range("Z1")
$tmp1.Font.Color = $black
$tmp1.Font.Size = 10
$tmp1.Font.Name = "Helvetica Neue"

$tmp2 = $ws.Range("Z2")
$tmp2.Font.Color = $black
$tmp2.Font.Size = 10
$tmp2.Font.Name = "Helvetica Neue"
$tmp2.Font.Bold = $true

$tmp1.Copy()
$ws.Range("N8:N57").PasteSpecial(-4122)
$ws.Range("A42:A47").PasteSpecial(-4122)

$tmp2.Copy()
$ws.Range("N7").PasteSpecial(-4122)

$ws.Range("Z1:Z2").Clear()
$excel.CutCopyMode = $false

Write-Output "done"
